$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 5 ("Timeline & Milestones"): bold the Phase 3 row in the table
# ---------------------------------------------------------------------
$sTimeline = $p.Slides.Item(5)
$tblTimeline = $sTimeline.Shapes.Item(3).Table
for ($c = 1; $c -le $tblTimeline.Columns.Count; $c++) {
    $tblTimeline.Cell(4, $c).Shape.TextFrame.TextRange.Font.Bold = -1
}

# ---------------------------------------------------------------------
# Slide 8 ("Investment Summary"): rework the investment table to add
# "AWS/Partner Credits" / net columns, and a new TOTAL INVESTMENT row
# ---------------------------------------------------------------------
$sInvest = $p.Slides.Item(8)
$tbl = $sInvest.Shapes.Item(3).Table

# Insert two new columns right after the existing "Year 1" column (col 2)
$tbl.Columns.Add(3) | Out-Null
$tbl.Columns.Add(3) | Out-Null

# Append a new row at the bottom for TOTAL INVESTMENT
$tbl.Rows.Add() | Out-Null

# Set the final column widths (EMU / 12700 = points)
$colWidths = @(1567967, 1132421, 1742186, 1132421, 993046, 993046, 1158554)
for ($c = 1; $c -le $colWidths.Length; $c++) {
    $tbl.Columns($c).Width = $colWidths[$c - 1] / 12700
}

# Ensure every row keeps the expected 370840 EMU (29.2pt) height
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $tbl.Rows($r).Height = 370840 / 12700
}

function Set-Cell($row, $col, $text, $bold, $size, $fillRgb, $fontRgb) {
    $cellShape = $tbl.Cell($row, $col).Shape
    $tr = $cellShape.TextFrame.TextRange
    $tr.Text = $text
    $tr.Font.Size = $size
    if ($bold) {
        $tr.Font.Bold = -1
    } else {
        $tr.Font.Bold = 0
    }
    if ($fontRgb -ne $null) {
        $tr.Font.Color.RGB = $fontRgb
    }
    $cellShape.Fill.ForeColor.RGB = $fillRgb
}

$headerFill = 138400     # A01C02 (BGR-packed)
$headerFont = 16777215   # FFFFFF
$dataFill   = 15132391   # E7E6E6

# ----- Header row -----
Set-Cell 1 1 "Cost Category" $true 14 $headerFill $headerFont
Set-Cell 1 2 "Year 1 List" $true 14 $headerFill $headerFont
Set-Cell 1 3 "AWS/Partner Credits" $true 14 $headerFill $headerFont
Set-Cell 1 4 "Year 1 Net" $true 14 $headerFill $headerFont
Set-Cell 1 5 "Year 2" $true 14 $headerFill $headerFont
Set-Cell 1 6 "Year 3" $true 14 $headerFill $headerFont
Set-Cell 1 7 "3-Year Total" $true 14 $headerFill $headerFont

# ----- Row 2: Professional Services -----
Set-Cell 2 1 "Professional Services" $false 11 $dataFill $null
Set-Cell 2 2 "`$82,250" $false 11 $dataFill $null
Set-Cell 2 3 "(`$10,000)" $false 11 $dataFill $null
Set-Cell 2 4 "`$72,250" $false 11 $dataFill $null
Set-Cell 2 5 "`$0" $false 11 $dataFill $null
Set-Cell 2 6 "`$0" $false 11 $dataFill $null
Set-Cell 2 7 "`$72,250" $false 11 $dataFill $null

# ----- Row 3: Cloud Infrastructure -----
Set-Cell 3 1 "Cloud Infrastructure" $false 11 $dataFill $null
Set-Cell 3 2 "`$26,830" $false 11 $dataFill $null
Set-Cell 3 3 "(`$5,000)" $false 11 $dataFill $null
Set-Cell 3 4 "`$21,830" $false 11 $dataFill $null
Set-Cell 3 5 "`$26,830" $false 11 $dataFill $null
Set-Cell 3 6 "`$26,830" $false 11 $dataFill $null
Set-Cell 3 7 "`$75,490" $false 11 $dataFill $null

# ----- Row 4: Software Licenses & Subscriptions -----
Set-Cell 4 1 "Software Licenses & Subscriptions" $false 11 $dataFill $null
Set-Cell 4 2 "`$2,784" $false 11 $dataFill $null
Set-Cell 4 3 "`$0" $false 11 $dataFill $null
Set-Cell 4 4 "`$2,784" $false 11 $dataFill $null
Set-Cell 4 5 "`$2,784" $false 11 $dataFill $null
Set-Cell 4 6 "`$2,784" $false 11 $dataFill $null
Set-Cell 4 7 "`$8,352" $false 11 $dataFill $null

# ----- Row 5: Support & Maintenance (was TOTAL SOLUTION INVESTMENT) -----
Set-Cell 5 1 "Support & Maintenance" $false 11 $dataFill $null
Set-Cell 5 2 "`$4,087" $false 11 $dataFill $null
Set-Cell 5 3 "`$0" $false 11 $dataFill $null
Set-Cell 5 4 "`$4,087" $false 11 $dataFill $null
Set-Cell 5 5 "`$4,087" $false 11 $dataFill $null
Set-Cell 5 6 "`$4,087" $false 11 $dataFill $null
Set-Cell 5 7 "`$12,261" $false 11 $dataFill $null

# ----- Row 6: TOTAL INVESTMENT (new row) -----
Set-Cell 6 1 "TOTAL INVESTMENT" $true 11 $dataFill $null
Set-Cell 6 2 "`$115,951" $true 11 $dataFill $null
Set-Cell 6 3 "(`$15,000)" $true 11 $dataFill $null
Set-Cell 6 4 "`$100,951" $true 11 $dataFill $null
Set-Cell 6 5 "`$33,701" $true 11 $dataFill $null
Set-Cell 6 6 "`$33,701" $true 11 $dataFill $null
Set-Cell 6 7 "`$168,353" $true 11 $dataFill $null
